$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update cell values per the diff (order matters for shared-string table order)
$ws.Range("A7").Value = "Unshielded VFD"
$ws.Range("A5").Value = "Crane with Unshielded VFD"
$ws.Range("A3").Value = "Microwave Dryer Take 2"
$ws.Range("A6").Value = "Microwave Dryer Take 1"
$ws.Range("F2").Value = "Hz"

# Update selection to F3
$ws.Range("F3").Select()
